# Generate Report for Handback
#
# For both locale sheets (zh-cn, de-de), row 6 corresponds to the
# "70e5338e-d868-485d-916f-e5d9708115d5" source file. A handback was
# produced for it, but it was not built against the latest handoff, so:
#   - "Latest Target File"      (col I) gets the handed-back source file
#                                 name, hyperlinked like col A/col I2.
#   - "Latest Handback File"    (col J) gets the generated xlf file name.
#   - "Latest Handback DateTime"(col K) gets the generation timestamp.
#   - "Error Detail"            (col P) gets the "not the latest" warning.
# Column P (Error Detail) is also widened to fit the longer message.

$wb = $excel.ActiveWorkbook

$handbackUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8d33de0daf95b15e8ad424f23985e4d82adf2d95/e2e/70e5338e-d868-485d-916f-e5d9708115d5.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2a44855bf84b2e004c1adcffacea34b2a1bf5081/e2e/70e5338e-d868-485d-916f-e5d9708115d5.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8d33de0daf95b15e8ad424f23985e4d82adf2d95/e2e/70e5338e-d868-485d-916f-e5d9708115d5.md."

$sheets = @{
    "zh-cn" = @{
        TargetFile   = "70e5338e-d868-485d-916f-e5d9708115d5.md"
        HandbackFile = "70e5338e-d868-485d-916f-e5d9708115d5.e95fec7d2c4eb53a44b214ac419766915dcb72f7.zh-cn.xlf"
        HandbackDate = "2016-08-22 00:54:29"
    }
    "de-de" = @{
        TargetFile   = "70e5338e-d868-485d-916f-e5d9708115d5.md"
        HandbackFile = "70e5338e-d868-485d-916f-e5d9708115d5.e95fec7d2c4eb53a44b214ac419766915dcb72f7.de-de.xlf"
        HandbackDate = "2016-08-22 00:54:36"
    }
}

foreach ($sheetName in $sheets.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $info = $sheets[$sheetName]

    # Latest Target File (I6) - file name + hyperlink to the handback commit,
    # styled the same way as the other filename hyperlinks on this sheet
    # (blue, underlined - matching the "HyperLink" cell style used in col A).
    $ws.Range("I6").Value = $info.TargetFile
    $ws.Hyperlinks.Add($ws.Range("I6"), $handbackUrl, "", "", $info.TargetFile)
    $ws.Range("I6").Font.Underline = 2
    $ws.Range("I6").Font.Color = 0xED9564
    $ws.Range("I6").Font.Name = "Calibri"

    # Latest Handback File (J6).
    $ws.Range("J6").Value = $info.HandbackFile

    # Latest Handback DateTime (K6).
    $ws.Range("K6").Value = $info.HandbackDate

    # Error Detail (P6).
    $ws.Range("P6").Value = $errorDetail

    # Widen the Error Detail column so the long message is readable.
    $ws.Columns.Item(16).ColumnWidth = 39.17
}
